$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated T_ic (column C) values (re-run of the Google Maps plotter calc) ---
$ws.Range("C2").Value  = 18.11364833999999
$ws.Range("C3").Value  = 21.415392116
$ws.Range("C4").Value  = 26.592446554999999
$ws.Range("C5").Value  = 33.036333196000001
$ws.Range("C6").Value  = 38.724609295500009
$ws.Range("C7").Value  = 43.946551491999998
$ws.Range("C8").Value  = 48.205900710000009
$ws.Range("C9").Value  = 14.198438577999999
$ws.Range("C10").Value = 10.9939516075
$ws.Range("C11").Value = 10.45779231
$ws.Range("C12").Value = 20.816888712499999
$ws.Range("C13").Value = 26.387957891999999

# --- Row 9 also got new Hic / Vic source readings ---
$ws.Range("A9").Value = 2264.375
$ws.Range("B9").Value = 250.91562500000001

# --- Selection moved from the T_ic column to the Hic/Vic columns ---
$ws.Range("A2:B13").Select()
